$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: client/project location -> new project path (was the BS project path, now the Givat Ram campus path)
$ws.Range("B3").Value = "W:\Projects\בהת\175 תכנית מתאר קמפוס גבעת רם\קבצי עבודה\תחזיות_דמוגרפיות"

# B4: scenario name (forecast_version) -> with_project
$ws.Range("B4").Value = "with_project"

# B5: v_date -> 241209
$ws.Range("B5").Value = 241209
